$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap row 5 and row 6 ---
$ws.Range("A5").Value2 = 130938725
$ws.Range("A6").Value2 = 130938736
$ws.Range("B5").Value2 = 91808
$ws.Range("B6").Value2 = 58043
$ws.Range("E5").Value2 = 1202
$ws.Range("E6").Value2 = 103021
$ws.Range("F5").Value2 = "Ullticka"
$ws.Range("F6").Value2 = "Talltita"
$ws.Range("G5").Value2 = "Phellinidium ferrugineofuscum"
$ws.Range("G6").Value2 = "Poecile montanus"
$ws.Range("H5").Value2 = "(P.Karst.) Fiasson & Niemelä"
$ws.Range("H6").Value2 = "(Conrad von Baldenstein, 1827)"
$ws.Range("I5").Value2 = ""
$ws.Range("I6").NumberFormat = "@"
$ws.Range("I6").Value2 = "1"
$ws.Range("J5").Value2 = ""
$ws.Range("J6").Value2 = ""
$ws.Range("K5").Value2 = "teleomorf"
$ws.Range("K6").Value2 = ""
$ws.Range("L5").Value2 = ""
$ws.Range("L6").Value2 = ""
$ws.Range("M5").Value2 = ""
$ws.Range("M6").Value2 = "förbiflygande"
$ws.Range("N5").Value2 = ""
$ws.Range("N6").Value2 = "observerad"
$ws.Range("Q5").Value2 = 476460
$ws.Range("Q6").Value2 = 476389
$ws.Range("R5").Value2 = 7033592
$ws.Range("R6").Value2 = 7033614
$ws.Range("AF5").Value2 = ""
$ws.Range("AF6").Value2 = ""
$ws.Range("AJ5").Value2 = "gran"
$ws.Range("AJ6").Value2 = ""
$ws.Range("AK5").Value2 = "Picea abies"
$ws.Range("AK6").Value2 = ""
$ws.Range("AM5").Value2 = "Liggande död trädstam, utan markontakt"
$ws.Range("AM6").Value2 = ""
$ws.Range("AO5").Value2 = "Horizontal, dead without ground contact # Picea abies"
$ws.Range("AO6").Value2 = ""

# --- Swap row 17 and row 18 ---
$ws.Range("A17").Value2 = 130938743
$ws.Range("A18").Value2 = 130938740
$ws.Range("B17").Value2 = 91828
$ws.Range("B18").Value2 = 78255
$ws.Range("E17").Value2 = 5432
$ws.Range("E18").Value2 = 228579
$ws.Range("F17").Value2 = "Granticka"
$ws.Range("F18").Value2 = "Liten svartspik"
$ws.Range("G17").Value2 = "Porodaedalea chrysoloma s.lat."
$ws.Range("G18").Value2 = "Chaenothecopsis nana"
$ws.Range("H17").Value2 = ""
$ws.Range("H18").Value2 = "Tibell"
$ws.Range("K17").Value2 = "teleomorf"
$ws.Range("K18").Value2 = ""
$ws.Range("Q17").Value2 = 476555
$ws.Range("Q18").Value2 = 476532
$ws.Range("R17").Value2 = 7033581
$ws.Range("R18").Value2 = 7033586
$ws.Range("AJ17").Value2 = "gran"
$ws.Range("AJ18").Value2 = ""
$ws.Range("AK17").Value2 = "Picea abies"
$ws.Range("AK18").Value2 = ""
$ws.Range("AO17").Value2 = "Picea abies"
$ws.Range("AO18").Value2 = ""

# --- Swap row 19 and row 20 ---
$ws.Range("A19").Value2 = 130938734
$ws.Range("A20").Value2 = 130938752
$ws.Range("B19").Value2 = 57884
$ws.Range("B20").Value2 = 79243
$ws.Range("E19").Value2 = 100109
$ws.Range("E20").Value2 = 6425
$ws.Range("F19").Value2 = "Tretåig hackspett"
$ws.Range("F20").Value2 = "Garnlav"
$ws.Range("G19").Value2 = "Picoides tridactylus"
$ws.Range("G20").Value2 = "Alectoria sarmentosa"
$ws.Range("H19").Value2 = "(Linnaeus, 1758)"
$ws.Range("H20").Value2 = "(Ach.) Ach."
$ws.Range("J19").Value2 = ""
$ws.Range("J20").Value2 = ""
$ws.Range("L19").Value2 = ""
$ws.Range("L20").Value2 = ""
$ws.Range("M19").Value2 = "äldre spår"
$ws.Range("M20").Value2 = ""
$ws.Range("Q19").Value2 = 476457
$ws.Range("Q20").Value2 = 476286
$ws.Range("R19").Value2 = 7033634
$ws.Range("R20").Value2 = 7033527
$ws.Range("AC19").Value2 = "Ringhack, äldre, på gran."
$ws.Range("AC20").Value2 = ""
$ws.Range("AF19").Value2 = ""
$ws.Range("AF20").Value2 = ""
$ws.Range("AM19").Value2 = "Trädstam på levande träd"
$ws.Range("AM20").Value2 = ""
$ws.Range("AO19").Value2 = "Stem on living tree # Picea abies"
$ws.Range("AO20").Value2 = "Picea abies"
